$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Insert a new row at position 9, shifting existing rows 9-24 down to 10-25.
$ws.Rows.Item(9).Insert()

# Row 9: brand-new job posting data.
$ws.Cells.Item(9, 1).Value = '2025-12-02 01:53:13'
$ws.Cells.Item(9, 2).Value = '【せどり×ツール製作】APIを使用したせどりツールを製作できるエンジニアさんを募集します♪'
$ws.Cells.Item(9, 3).Value = 'システム開発'
$ws.Cells.Item(9, 4).Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Cells.Item(9, 5).Value = '期限情報なし'
$ws.Cells.Item(9, 6).Value = 'https://www.lancers.jp/work/detail/5217096'
$ws.Cells.Item(9, 7).Value = 243
$ws.Cells.Item(9, 8).Value = '🔥API ◆ツール'

# Rows 2-8: refresh retrieval timestamp only.
$newTimestamp = '2025-12-02 01:53:13'
foreach ($r in 2..8) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# Rows 10-25: content already shifted down by the row insert; refresh timestamp only.
foreach ($r in 10..25) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# New hyperlink relationship for the freshly appended row 25 (matches source export behavior).
$ws.Hyperlinks.Add($ws.Range("F25"), "https://www.lancers.jp/work/detail/5444370") | Out-Null

$ws.Range("A1").Select()
